$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold/border/center style) from existing header cell O1 to the new header cells P1, Q1
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))

# Full target grid: rows 1-25 (index 0-24), columns A-Q (index 0-16); $null = cell not present in sheet
$data = @(
        @($null,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15),
        @(0,3.516539377158892,1.078588852915516,0.04350420173832958,1.419166398039707,0.472515791760344,0.0007792636016987003,0,0.009407384440288435,0,0,0,0,0,0,0,1.277260632581914),
        @(1,3.060779801384797,0.9511387153850706,0.03898626743634281,1.218044655174424,0.4310114202132667,0.0007839042084106291,0,0.006021793909746886,0,0,0,0,0,0,0,1.189011944184756),
        @(2,2.78116008243444,0.87359336389639,0.03623777507969805,1.096911651310066,0.4066227789309877,0.0007868358479447618,0,0.004317012687420352,0,0,0,0,0,0,0,1.138426943591753),
        @(3,2.667191196798228,0.8437270183374039,0.03520616445518954,1.048043604885592,0.3961481752877631,0.0007880570292950782,0,0.00376135697338098,0,0,0,0,0,0,0,1.116017533858042),
        @(4,2.64822199796032,0.8407106350256299,0.03513649939789332,1.039947653828108,0.3934623566796347,0.0007882673851986599,0,0.003746778042270194,0,0,0,0,0,0,0,1.109142225967361),
        @(5,2.779499325879272,0.8784884188306989,0.03650150033514166,1.096224134122892,0.4038525123599754,0.0007868689694541176,0,0.004491715281754516,0,0,0,0,0,0,0,1.129367300635408),
        @(6,3.35915009768047,1.041581425468934,0.04231366816730997,1.349232821394779,0.4544744364173638,0.0007808674799953397,0,0.008365928385478227,0,0,0,0,0,0,0,1.234434928888305),
        @(7,4.499849896746582,1.35844943808803,0.05343420173026203,1.868677731637362,0.5683922495544991,0.0007696917840793197,0,0.01908626117053647,0,0,0,0,0,0,0,1.488660917460436),
        @(8,5.343344342455623,1.599679545547531,0.0619309849569305,2.271600381442127,0.6587116949687726,0.0007618150989550717,0,0.0297290123805638,0,0,0,0,0,0,0,1.696892671837361),
        @(9,5.728812441600951,1.720208489292531,0.06633152740253223,2.461300675634007,0.6976538564766344,0.0007583098459426427,0,0.03556790276178123,0,0,0,0,0,0,0,1.784373770071852),
        @(10,5.875273707070789,1.761292681960299,0.06775030166279095,2.534289032760199,0.7150632417429961,0.0007569789155685768,0,0.0377610249083471,0,0,0,0,0,0,0,1.826287979662226),
        @(11,5.843740394889835,1.75145372885612,0.06739250485865966,2.518525732698905,0.7117606621950046,0.0007572631490899412,0,0.03725280503380191,0,0,0,0,0,0,0,1.818739865924272),
        @(12,5.740865619365707,1.72316769101127,0.06642607279552237,2.467285667497748,0.6992754886482118,0.0007581993115455653,0,0.03573261709733533,0,0,0,0,0,0,0,1.788448826704069),
        @(13,5.677848752389139,1.707801802638642,0.06593716203332178,2.436033870324636,0.690764591165717,0.0007587778398641799,0,0.03488187362864537,0,0,0,0,0,0,0,1.767035982827963),
        @(14,5.317763693821632,1.607965184848581,0.06249831337174783,2.259216006557381,0.6485290530646637,0.0007620815860210559,0,0.02985493460792554,0,0,0,0,0,0,0,1.665891769812589),
        @(15,5.097431413053698,1.547085964193855,0.06039792244490627,2.152528401578252,0.6232469611579887,0.0007641225321531033,0,0.02697771689609052,0,0,0,0,0,0,0,1.605876999379745),
        @(16,4.971015782787504,1.507314701841892,0.05893465276424337,2.091801238326426,0.611239715954099,0.0007652907878438075,0,0.02520486716114156,0,0,0,0,0,0,0,1.579729249949793),
        @(17,4.928144846896544,1.497632434803165,0.05863895716920808,2.071309092459586,0.6054217223431877,0.0007656959796628295,0,0.02475889296226175,0,0,0,0,0,0,0,1.565074211779915),
        @(18,5.120876494190384,1.553169766999133,0.06060085227606038,2.163825200279049,0.6260974780276314,0.0007639038309240136,0,0.02726146111278283,0,0,0,0,0,0,0,1.612793806254501),
        @(19,5.770921346933903,1.736956867280469,0.067000802567172,2.482256383194965,0.7003351339900092,0.0007579359506373018,0,0.03635845982994201,0,0,0,0,0,0,0,1.78872575005883),
        @(20,6.198113164505116,1.851440491424967,0.07084702028519274,2.69702540378313,0.7543929982663542,0.0007540623535256584,0,0.04280160602805338,0,0,0,0,0,0,0,1.921863553715895),
        @(21,5.970073204496771,1.783587723937899,0.06844000897443436,2.581788937061958,0.7284550675454113,0.0007561122197995173,0,0.0390887220535463,0,0,0,0,0,0,0,1.860455886910245),
        @(22,5.110550384075509,1.540093550188203,0.05996250411526916,2.158797249292547,0.6297227669184764,0.0007639790849492597,0,0.02677466222253155,0,0,0,0,0,0,0,1.625930566867396),
        @(23,4.190293018286468,1.281730338667899,0.0509179980049197,1.724912907644182,0.5315540212227958,0.0007726787525529291,0,0.01606038940327714,0,0,0,0,0,0,0,1.399985515818884)
)

for ($r = 1; $r -le 25; $r++) {
    for ($c = 1; $c -le 17; $c++) {
        $v = $data[$r - 1][$c - 1]
        if ($v -ne $null) {
            $ws.Cells.Item($r, $c).Value = $v
        }
    }
}